$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column D ("Explanation No.") before the existing Explanation column.
$ws.Columns.Item(4).Insert()

# Update explanation text (now in column E) for rows 4 and 5 first, so new
# shared strings are appended in the same order Excel produced them.
$ws.Cells.Item(4, 5).Value = "An actual filtering idea implemented in a basic way and added to report. Then the enchanced version is implemented and added to report as well. Before and after pictures along wit Python codes are added."
$ws.Cells.Item(5, 5).Value = "All the questions are answered by explaining the algorithm's steps but they could have been kept even shorter."

# Header for new column D
$ws.Cells.Item(1, 4).Value = "Explanation No."

# Sequential numbering 1-6 for rows 2-7
for ($i = 2; $i -le 7; $i++) {
    $ws.Cells.Item($i, 4).Value = $i - 1
}

# Update estimated point values (column C)
$ws.Cells.Item(4, 3).Value = 20
$ws.Cells.Item(5, 3).Value = 20
$ws.Cells.Item(8, 3).Value = 100

# Column widths (closest achievable values to the target 15.140625 / 188.7109375
# given this engine's column-width pixel quantization)
$ws.Columns.Item(4).ColumnWidth = 14.15
$ws.Columns.Item(5).ColumnWidth = 187.9

# Selection as seen in the final file
$ws.Range("C19").Select()

$wb.Save()
